# envio preliminar IRA 42 2025
# Adds week-42 column (AS) to the weekly IRA/UCI report sheet and
# updates a few institution names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for week 42 (column AS), styled like the other header cells ---
# Force text (like the other week-number headers) instead of Excel's
# automatic numeric conversion of a digit-only value.
$ws.Range("AS1").NumberFormat = "@"
$ws.Range("AS1").Value = "42"
$ws.Range("AS1").Font.Bold = $true
$ws.Range("AS1").HorizontalAlignment = -4108   # xlCenter

# --- Institution name corrections (column C) ---
$ws.Range("C5").Value  = "CAJA DE COMPENSACION FAMILIAR DE RISARALDA COMFAMI"
$ws.Range("C6").Value  = "CAJA DE COMPENSACION FAMILIAR DE RISARALDA COMFAMI"
$ws.Range("C44").Value = "EMPRESA DE MEDICINA INTEGRAL EMI SA - SERVICIO DE"

# --- Week 42 data values for column AS (only rows with reported data) ---
$as_values = @{
    2  = 0
    3  = 0
    5  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    14 = 0
    16 = 0
    17 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    28 = 1
    29 = 1
    30 = 2
    31 = 0
    36 = 0
    37 = 0
    38 = 0
    41 = 0
    42 = 0
    43 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 0
    57 = 0
    58 = 0
}

# Row 28 previously had no value in column AR either; fill it in now.
$ws.Range("AR28").Value = 0

foreach ($row in $as_values.Keys) {
    $ws.Range("AS$row").Value = $as_values[$row]
}
